$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1587539.5
$ws.Range("I6").Value = 1904847.4
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 5714542.199999999
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -5714430.199999999
$ws.Range("N6").Value = -3224
$ws.Range("H129").Value = 1971.6316
$ws.Range("I129").Value = 1139.25
$ws.Range("J129").Value = 2193.6
$ws.Range("K129").Value = 3417.75
$ws.Range("L129").Value = 6580.799999999999
$ws.Range("M129").Value = 1582.25
$ws.Range("N129").Value = -16580.8
$ws.Range("H132").Value = 6061.673
$ws.Range("I132").Value = 6838.6665
$ws.Range("J132").Value = 5312.4287
$ws.Range("K132").Value = 20515.9995
$ws.Range("L132").Value = 15937.2861
$ws.Range("M132").Value = -17985.9995
$ws.Range("N132").Value = -20997.2861
$ws.Range("H135").Value = 594.38
$ws.Range("I135").Value = 250.1282
$ws.Range("J135").Value = 1814.909
$ws.Range("K135").Value = 2251.1538
$ws.Range("L135").Value = 16334.181
$ws.Range("M135").Value = 283.8462
$ws.Range("N135").Value = -21404.181
$ws.Range("H137").Value = 1303.3693
$ws.Range("I137").Value = 1554.6389
$ws.Range("J137").Value = 991.4483
$ws.Range("K137").Value = 4663.9167
$ws.Range("L137").Value = 2974.3449
$ws.Range("M137").Value = -2113.9167
$ws.Range("N137").Value = -8074.3449
$ws.Range("H138").Value = 1606.2373
$ws.Range("I138").Value = 1053.8182
$ws.Range("J138").Value = 3226.6667
$ws.Range("K138").Value = 3161.4546
$ws.Range("L138").Value = 9680.000100000001
$ws.Range("M138").Value = 1978.5454
$ws.Range("N138").Value = -19960.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 4000031
$ws.Range("I11").Value = 4000031
$ws.Range("K11").Value = 4000031
$ws.Range("M11").Value = -3999887
$ws.Range("H63").Value = 3628.5
$ws.Range("I63").Value = 2418.9375
$ws.Range("J63").Value = 8466.75
$ws.Range("K63").Value = 2418.9375
$ws.Range("L63").Value = 8466.75
$ws.Range("M63").Value = -1732.9375
$ws.Range("N63").Value = -9838.75
$ws.Range("H66").Value = 3628.5
$ws.Range("I66").Value = 2418.9375
$ws.Range("J66").Value = 8466.75
$ws.Range("K66").Value = 12094.6875
$ws.Range("L66").Value = 42333.75
$ws.Range("M66").Value = -8662.6875
$ws.Range("N66").Value = -49197.75
$ws.Range("H74").Value = 1414.5555
$ws.Range("I74").Value = 1044.4474
$ws.Range("J74").Value = 2293.5625
$ws.Range("K74").Value = 1044.4474
$ws.Range("L74").Value = 2293.5625
$ws.Range("M74").Value = -170.4474
$ws.Range("N74").Value = -4041.5625
$ws.Range("H77").Value = 1414.5555
$ws.Range("I77").Value = 1044.4474
$ws.Range("J77").Value = 2293.5625
$ws.Range("K77").Value = 5222.237
$ws.Range("L77").Value = 11467.8125
$ws.Range("M77").Value = -854.2370000000001
$ws.Range("N77").Value = -20203.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4732.9624
$ws.Range("I134").Value = 2126.261
$ws.Range("J134").Value = 6731.433
$ws.Range("K134").Value = 6378.782999999999
$ws.Range("L134").Value = 20194.299
$ws.Range("M134").Value = -3843.782999999999
$ws.Range("N134").Value = -25264.299

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1782.875
$ws.Range("I10").Value = 374.5
$ws.Range("J10").Value = 6008
$ws.Range("K10").Value = 374.5
$ws.Range("L10").Value = 6008
$ws.Range("M10").Value = -235.5
$ws.Range("N10").Value = -6286
$ws.Range("H31").Value = 5052208.5
$ws.Range("I31").Value = 1408.6327
$ws.Range("J31").Value = 19610398
$ws.Range("K31").Value = 1408.6327
$ws.Range("L31").Value = 19610398
$ws.Range("M31").Value = -1113.6327
$ws.Range("N31").Value = -19610988
$ws.Range("H34").Value = 5052208.5
$ws.Range("I34").Value = 1408.6327
$ws.Range("J34").Value = 19610398
$ws.Range("K34").Value = 1408.6327
$ws.Range("L34").Value = 19610398
$ws.Range("M34").Value = -1206.6327
$ws.Range("N34").Value = -19610802
$ws.Range("H58").Value = 1639.5
$ws.Range("I58").Value = 1249.9863
$ws.Range("J58").Value = 2692.6296
$ws.Range("K58").Value = 1249.9863
$ws.Range("L58").Value = 2692.6296
$ws.Range("M58").Value = -1046.9863
$ws.Range("N58").Value = -3098.6296
$ws.Range("H136").Value = 1639.5
$ws.Range("I136").Value = 1249.9863
$ws.Range("J136").Value = 2692.6296
$ws.Range("K136").Value = 3749.9589
$ws.Range("L136").Value = 8077.888800000001
$ws.Range("M136").Value = -1199.9589
$ws.Range("N136").Value = -13177.8888

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 3377.4194
$ws.Range("I23").Value = 43.666668
$ws.Range("J23").Value = 4177.52
$ws.Range("K23").Value = 131.000004
$ws.Range("L23").Value = 12532.56
$ws.Range("M23").Value = 103.999996
$ws.Range("N23").Value = -13002.56
$ws.Range("H33").Value = 170.44444
$ws.Range("I33").Value = 66
$ws.Range("J33").Value = 254
$ws.Range("K33").Value = 396
$ws.Range("L33").Value = 1524
$ws.Range("M33").Value = -113
$ws.Range("N33").Value = -2090
$ws.Range("H122").Value = 2392.1228
$ws.Range("I122").Value = 982.9091
$ws.Range("J122").Value = 2729.1086
$ws.Range("K122").Value = 8846.1819
$ws.Range("L122").Value = 24561.9774
$ws.Range("M122").Value = -6396.1819
$ws.Range("N122").Value = -29461.9774

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 12195.1875
$ws.Range("I102").Value = 13432.889
$ws.Range("J102").Value = 10603.857
$ws.Range("K102").Value = 13432.889
$ws.Range("L102").Value = 10603.857
$ws.Range("M102").Value = -11810.889
$ws.Range("N102").Value = -13847.857
$ws.Range("H107").Value = 803.0952
$ws.Range("I107").Value = 719.9286
$ws.Range("J107").Value = 969.4286
$ws.Range("K107").Value = 719.9286
$ws.Range("L107").Value = 969.4286
$ws.Range("M107").Value = 1200.0714
$ws.Range("N107").Value = -4809.4286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2600.3333
$ws.Range("I61").Value = 2320.5
$ws.Range("K61").Value = 2320.5
$ws.Range("M61").Value = -2118.5
$ws.Range("H113").Value = 2600.3333
$ws.Range("I113").Value = 2320.5
$ws.Range("K113").Value = 2320.5
$ws.Range("M113").Value = -150.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3426.3845
$ws.Range("I126").Value = 4200.3
$ws.Range("J126").Value = 846.6667
$ws.Range("K126").Value = 12600.9
$ws.Range("L126").Value = 2540.0001
$ws.Range("M126").Value = -10130.9
$ws.Range("N126").Value = -7480.0001
$ws.Range("H132").Value = 1319.2638
$ws.Range("I132").Value = 1123.4445
$ws.Range("J132").Value = 1605.0541
$ws.Range("K132").Value = 3370.3335
$ws.Range("L132").Value = 4815.1623
$ws.Range("M132").Value = -840.3335000000002
$ws.Range("N132").Value = -9875.1623
